$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$bio = "Sadie is a Research Associate I for the McCarroll & Macosko BICAN project. She is a recent graduate from Wheaton College MA with a B.S. in Neuroscience on the Pre-Medical Track. Prior to joining the BICAN project, in her Behavioral Neuroscience lab she completed a senior honors thesis titled: The effects of maternal separation and social isolation on memory and myelin in adolescent rats. In addition to being a research associate, she is also a medical assistant in primary care and sports medicine at DMC primary care. "
$email = "drouin@broadinstitute.org"

# --- Add new row 36 to Sheet1 for the newest team member ---
$ws1.Cells.Item(36, 1).Value = "Sadie"
$ws1.Cells.Item(36, 2).Value = "Drouin"
$ws1.Cells.Item(36, 7).Value = $bio
$ws1.Cells.Item(36, 6).Value = $email
$ws1.Cells.Item(36, 8).Value = "Data Generation"

$ws1.Hyperlinks.Add($ws1.Cells.Item(36, 6), "mailto:$email") | Out-Null
$ws1.Cells.Item(36, 6).Style = "Hyperlink"

# --- Replace the Sheet2 staging row with the new person's info ---
$ws2.Cells.Item(2, 1).Value = "Sadie"
$ws2.Cells.Item(2, 2).Value = "Drouin"
$ws2.Cells.Item(2, 3).ClearContents()
$ws2.Cells.Item(2, 4).ClearContents()
$ws2.Cells.Item(2, 6).Value = $bio
$ws2.Cells.Item(2, 7).Value = $email
$ws2.Cells.Item(2, 8).Value = "Data Generation"

$ws2.Hyperlinks.Add($ws2.Cells.Item(2, 6), "mailto:$email") | Out-Null
$ws2.Cells.Item(2, 6).Style = "Hyperlink"

# --- Select the newly added row on Sheet1, then switch to Sheet2 ---
$ws1.Range("A36:H36").Select() | Out-Null
$ws2.Activate()
